# Sincronização de dados: insere um novo registro de avaliação na linha 3
# (empurrando as linhas existentes 3..29 para 4..30) e preenche seus valores.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insere uma nova linha vazia em 3, deslocando as linhas abaixo para baixo.
$ws.Rows(3).Insert()

# Preenche os dados da nova avaliação (stars, comment, createdAt, task.id).
$ws.Range("A3").Value = 4
$ws.Range("C3").Value = 46001.65142589121
$ws.Range("D3").Value = "YTU1MjU5YTEtYTBkNy00ZDg5LWE5ZDAtN2M1MDIyODc2ZWZhOjU3MDE2"
